$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 25313.889
$ws.Range("I97").Value = 2916.6667
$ws.Range("J97").Value = 36512.5
$ws.Range("K97").Value = 8750.000100000001
$ws.Range("L97").Value = 109537.5
$ws.Range("M97").Value = -8254.000100000001
$ws.Range("N97").Value = -110529.5
$ws.Range("H99").Value = 265.42856
$ws.Range("I99").Value = 265.42856
$ws.Range("K99").Value = 796.28568
$ws.Range("M99").Value = 701.71432
$ws.Range("H103").Value = 445.33334
$ws.Range("J103").Value = 349.25
$ws.Range("L103").Value = 1047.75
$ws.Range("N103").Value = -2219.75
$ws.Range("H121").Value = 2287.8
$ws.Range("J121").Value = 2287.8
$ws.Range("L121").Value = 6863.400000000001
$ws.Range("N121").Value = -10357.4
$ws.Range("H125").Value = 2271.4
$ws.Range("I125").Value = 2162
$ws.Range("J125").Value = 2435.5
$ws.Range("K125").Value = 19458
$ws.Range("L125").Value = 21919.5
$ws.Range("M125").Value = -16998
$ws.Range("N125").Value = -26839.5
$ws.Range("H131").Value = 7746.4644
$ws.Range("J131").Value = 27229.428
$ws.Range("L131").Value = 81688.284
$ws.Range("N131").Value = -91768.284
$ws.Range("H132").Value = 73707.92999999999
$ws.Range("I132").Value = 1875.625
$ws.Range("J132").Value = 169484.33
$ws.Range("K132").Value = 5626.875
$ws.Range("L132").Value = 508452.99
$ws.Range("M132").Value = -3096.875
$ws.Range("N132").Value = -513512.99
$ws.Range("H138").Value = 6175404
$ws.Range("I138").Value = 1221.7142
$ws.Range("J138").Value = 9437236
$ws.Range("K138").Value = 3665.1426
$ws.Range("L138").Value = 28311708
$ws.Range("M138").Value = 1474.8574
$ws.Range("N138").Value = -28321988

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 587
$ws.Range("I4").Value = 211.75
$ws.Range("J4").Value = 1337.5
$ws.Range("K4").Value = 211.75
$ws.Range("L4").Value = 1337.5
$ws.Range("M4").Value = -95.75
$ws.Range("N4").Value = -1569.5
$ws.Range("H61").Value = 3711.2856
$ws.Range("I61").Value = 3313.739
$ws.Range("K61").Value = 3313.739
$ws.Range("M61").Value = -3101.739
$ws.Range("H122").Value = 1575.2307
$ws.Range("I122").Value = 1225.5294
$ws.Range("K122").Value = 3676.5882
$ws.Range("M122").Value = -1226.5882
$ws.Range("H132").Value = 2297.9792
$ws.Range("I132").Value = 2047.8334
$ws.Range("K132").Value = 6143.5002
$ws.Range("M132").Value = -3613.5002
$ws.Range("H136").Value = 3711.2856
$ws.Range("I136").Value = 3313.739
$ws.Range("K136").Value = 9941.217000000001
$ws.Range("M136").Value = -7391.217000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4138.15
$ws.Range("I99").Value = 1978.5
$ws.Range("K99").Value = 1978.5
$ws.Range("M99").Value = -480.5
$ws.Range("H105").Value = 2049.0667
$ws.Range("I105").Value = 2059.7144
$ws.Range("K105").Value = 2059.7144
$ws.Range("M105").Value = -312.7143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 74163.07000000001
$ws.Range("I31").Value = 85376.164
$ws.Range("J31").Value = 6884.5
$ws.Range("K31").Value = 85376.164
$ws.Range("L31").Value = 6884.5
$ws.Range("M31").Value = -85081.164
$ws.Range("N31").Value = -7474.5
$ws.Range("H34").Value = 74163.07000000001
$ws.Range("I34").Value = 85376.164
$ws.Range("J34").Value = 6884.5
$ws.Range("K34").Value = 85376.164
$ws.Range("L34").Value = 6884.5
$ws.Range("M34").Value = -85174.164
$ws.Range("N34").Value = -7288.5
$ws.Range("H58").Value = 3471.8276
$ws.Range("I58").Value = 3215.1
$ws.Range("K58").Value = 3215.1
$ws.Range("M58").Value = -3012.1
$ws.Range("H60").Value = 18937.375
$ws.Range("J60").Value = 18785.715
$ws.Range("L60").Value = 18785.715
$ws.Range("N60").Value = -19807.715
$ws.Range("H62").Value = 7610.5557
$ws.Range("I62").Value = 6217.8
$ws.Range("K62").Value = 6217.8
$ws.Range("M62").Value = -5593.8
$ws.Range("H65").Value = 7610.5557
$ws.Range("I65").Value = 6217.8
$ws.Range("K65").Value = 31089
$ws.Range("M65").Value = -27969
$ws.Range("H132").Value = 3534.55
$ws.Range("I132").Value = 3497.7856
$ws.Range("K132").Value = 10493.3568
$ws.Range("M132").Value = -7963.356800000001
$ws.Range("H134").Value = 18812.074
$ws.Range("I134").Value = 10967.904
$ws.Range("J134").Value = 46266.668
$ws.Range("K134").Value = 32903.712
$ws.Range("L134").Value = 138800.004
$ws.Range("M134").Value = -30368.712
$ws.Range("N134").Value = -143870.004
$ws.Range("H136").Value = 3471.8276
$ws.Range("I136").Value = 3215.1
$ws.Range("K136").Value = 9645.299999999999
$ws.Range("M136").Value = -7095.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 1500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4388
$ws.Range("N5").ClearContents()
$ws.Range("H113").Value = 458.57144
$ws.Range("I113").Value = 297.25
$ws.Range("K113").Value = 891.75
$ws.Range("M113").Value = 1278.25
$ws.Range("H121").Value = 781.6667
$ws.Range("J121").Value = 932
$ws.Range("L121").Value = 2796
$ws.Range("N121").Value = -5416
$ws.Range("H135").Value = 1500
$ws.Range("I135").Value = 1500
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13500
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -10965
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 3141
$ws.Range("I140").Value = 3141
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 9423
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -4243
$ws.Range("N140").Value = -22258.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2233.8667
$ws.Range("I80").Value = 2458.4285
$ws.Range("J80").Value = 2037.375
$ws.Range("K80").Value = 2458.4285
$ws.Range("L80").Value = 2037.375
$ws.Range("M80").Value = -1460.4285
$ws.Range("N80").Value = -4033.375
$ws.Range("H83").Value = 2233.8667
$ws.Range("I83").Value = 2458.4285
$ws.Range("J83").Value = 2037.375
$ws.Range("K83").Value = 12292.1425
$ws.Range("L83").Value = 10186.875
$ws.Range("M83").Value = -7300.1425
$ws.Range("N83").Value = -20170.875
$ws.Range("H132").Value = 3206.8635
$ws.Range("I132").Value = 3033.5625
$ws.Range("J132").Value = 3669
$ws.Range("K132").Value = 9100.6875
$ws.Range("L132").Value = 11007
$ws.Range("M132").Value = -6570.6875
$ws.Range("N132").Value = -16067

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 4999.5
$ws.Range("I8").Value = 4999
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 4999
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -4859
$ws.Range("N8").Value = -5280
$ws.Range("H122").Value = 2897.5334
$ws.Range("I122").Value = 2527.5557
$ws.Range("K122").Value = 7582.6671
$ws.Range("M122").Value = -5132.6671
$ws.Range("H131").Value = 99499.75
$ws.Range("J131").Value = 99499.75
$ws.Range("L131").Value = 99499.75
$ws.Range("N131").Value = -109579.75
$ws.Range("H132").Value = 1907.4117
$ws.Range("I132").Value = 1930.25
$ws.Range("J132").Value = 1800.8334
$ws.Range("K132").Value = 5790.75
$ws.Range("L132").Value = 5402.5002
$ws.Range("M132").Value = -3260.75
$ws.Range("N132").Value = -10462.5002
